$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing task descriptions (column A) ---
$ws.Range("A2").Value = "Setting up project - double jump, adding rigidbody etc"
$ws.Range("A3").Value = "Writing code for double jump "
$ws.Range("A5").Value = "Setting up project - movement and camera follow"
$ws.Range("A6").Value = "Writing camera follow code "
$ws.Range("A7").Value = "Writing Player Movement code "

# --- Add new rows 9-12 for the 02/03/2021 (pick-up script) session ---
$ws.Range("A9").Value = "Setting project - pick up "
$ws.Range("B9").Value = 44257
$ws.Range("C9").Value = 0.05
$ws.Range("D9").Value = 0.54861111111111105
$ws.Range("E9").Value = 0.55069444444444449
$ws.Range("G9").Value = 0.03

$ws.Range("A10").Value = "adding rigidbody to player and adding movement code "
$ws.Range("B10").Value = 44257
$ws.Range("C10").Value = 0.3
$ws.Range("D10").Value = 0.55138888888888882
$ws.Range("E10").Value = 0.57500000000000007
$ws.Range("F10").Value = 0.04
$ws.Range("G10").Value = 0.34

$ws.Range("A11").Value = "fixing broken script "
$ws.Range("B11").Value = 44257
$ws.Range("C11").Value = 0.15
$ws.Range("D11").Value = 0.57638888888888895
$ws.Range("E11").Value = 0.59722222222222221
$ws.Range("G11").Value = 0.3

$ws.Range("A12").Value = "Pick Up script "
$ws.Range("B12").Value = 44257
$ws.Range("C12").Value = 0.35
$ws.Range("D12").Value = 0.60069444444444442
$ws.Range("E12").Value = 0.6333333333333333
$ws.Range("F12").Value = 0.12
$ws.Range("G12").Value = 0.47

# Reuse the same date/time number formatting already used by row 7 so no
# redundant style entries get created.
$ws.Range("B7").Copy()
$ws.Range("B9:B12").PasteSpecial(-4122)
$ws.Range("D7:E7").Copy()
$ws.Range("D9:E12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Column A widened to fit the longer task descriptions ---
$ws.Columns.Item(1).ColumnWidth = 49.5

# --- Selection moves to A13 (first empty row) ---
$ws.Range("A13").Select() | Out-Null
